Write-Output ($ppt | Get-Member)
